$d = $word.ActiveDocument

# Update author first name: "Jane" -> "Hayden"
$d.Content.Find.Execute("Jane", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Hayden", 2)

# Update author last name: "Doe" -> "Atchley"
$d.Content.Find.Execute("Doe", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Atchley", 2)

# Update date: "1/31/23" -> "31 January 2023"
$d.Content.Find.Execute("1/31/23", $false, $false, $false, $false, $false,
                         $true, 1, $false, "31 January 2023", 2)
